# Batch solving and output ignition and solver times to input file;
# keyboardInterrupt NOT accounted for

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VarPropEta")

# --- Row 31: add Flux [MW/m^2] (F31) and t_ign [ms] (J31) values ---
$ws.Range("F31").Value = 0.8
$ws.Range("J31").Value = 0.27

# --- Row 32: fill in a new "Conserv-5" data row ---
$ws.Range("A32").Value = "Conserv-5"
$ws.Range("B32").Value = 101
$ws.Range("C32").Value = 601
$ws.Range("D32").Value = 400
$ws.Range("E32").Value = 800
$ws.Range("F32").Value = "n/a"
$ws.Range("G32").Value = 70
$ws.Range("H32").Value = 4890000
$ws.Range("I32").Value = 63
$ws.Range("J32").Value = 0.096

# --- Row 33: fill in a new "Conserv-6" data row ---
$ws.Range("A33").Value = "Conserv-6"
$ws.Range("B33").Value = 101
$ws.Range("C33").Value = 601
$ws.Range("D33").Value = 400
$ws.Range("E33").Value = 800
$ws.Range("F33").Value = "n/a"
$ws.Range("G33").Value = 48
$ws.Range("H33").Value = 4890000
$ws.Range("I33").Value = 63
$ws.Range("J33").Value = 0.045

# --- Row 34: fill in values (A34 stays blank) ---
$ws.Range("B34").Value = 101
$ws.Range("C34").Value = 601
$ws.Range("D34").Value = 400
$ws.Range("E34").Value = 1200
$ws.Range("F34").Value = "n/a"
$ws.Range("G34").Value = 48
$ws.Range("H34").Value = 4890000
$ws.Range("I34").Value = 63

# --- Row 35: fill in values (A35 stays blank) ---
$ws.Range("B35").Value = 101
$ws.Range("C35").Value = 601
$ws.Range("D35").Value = 200
$ws.Range("E35").Value = 1200
$ws.Range("F35").Value = "n/a"
$ws.Range("G35").Value = 48
$ws.Range("H35").Value = 4890000
$ws.Range("I35").Value = 63

# --- Style the F column so it matches the "n/a" formatted cells already ---
$ws.Range("F34").Style = $ws.Range("F10").Style
$ws.Range("F35").Style = $ws.Range("F10").Style

# --- View state: frozen pane / selection moved while scrolling through the new rows ---
$ws.Activate()
$ws.Range("A27").Select()
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("K31").Select()
